# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Primera/Segunda) for "Pepino ensalada"
# at sheet rows 76-77, shifting the existing data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 76; Excel shifts rows 76:187 down to 78:189
$ws.Rows("76:77").Insert()

# New row 76: "Primera" quality entry for date 2021-10-08 (serial 44477)
$ws.Cells.Item(76, 1).Value = 1
$ws.Cells.Item(76, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(76, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(76, 4).Value = 44477
$ws.Cells.Item(76, 5).Value = 15
$ws.Cells.Item(76, 6).Value = 100112043
$ws.Cells.Item(76, 7).Value = "Pepino ensalada"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 120
$ws.Cells.Item(76, 11).Value = 14000
$ws.Cells.Item(76, 12).Value = 15000
$ws.Cells.Item(76, 13).Value = 14500
$ws.Cells.Item(76, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(76, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(76, 16).Value = 207
$ws.Cells.Item(76, 17).Value = 70
$ws.Cells.Item(76, 18).Value = "Hortaliza"

# New row 77: "Segunda" quality entry for date 2021-10-08 (serial 44477)
$ws.Cells.Item(77, 1).Value = 1
$ws.Cells.Item(77, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(77, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(77, 4).Value = 44477
$ws.Cells.Item(77, 5).Value = 15
$ws.Cells.Item(77, 6).Value = 100112043
$ws.Cells.Item(77, 7).Value = "Pepino ensalada"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Segunda"
$ws.Cells.Item(77, 10).Value = 140
$ws.Cells.Item(77, 11).Value = 10000
$ws.Cells.Item(77, 12).Value = 11000
$ws.Cells.Item(77, 13).Value = 10500
$ws.Cells.Item(77, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(77, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(77, 16).Value = 105
$ws.Cells.Item(77, 17).Value = 100
$ws.Cells.Item(77, 18).Value = "Hortaliza"

# Ensure the date cells keep the date-formatted style (s="2" like column D elsewhere)
$ws.Cells.Item(76, 4).NumberFormat = $ws.Cells.Item(78, 4).NumberFormat
$ws.Cells.Item(77, 4).NumberFormat = $ws.Cells.Item(78, 4).NumberFormat
